$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")
$ws.Activate()

# Update data values on row 3
$ws.Range("G3").Value = 30697657335
$ws.Range("M3").Value = 307

# Update the view: scroll so D1 becomes the top-left visible cell, then
# move the selection to M4 (a single cell) to match the new saved view.
$excel.ActiveWindow.ScrollColumn = 4
$excel.ActiveWindow.ScrollRow = 1
$ws.Range("M4").Select()
